# Fix Training Data Issue (#48)
# Data in column BF ("Date") was stored as "12-24-2007-08" (a mangled
# concatenation of the game date and the season label). Correct it to the
# actual ISO-style game date "2007-12-24" for every data row.
#
# The BF column values must remain plain text (not be re-interpreted by
# Excel as a date serial number), so the cells are formatted as Text
# ("@") before the corrected string is written into them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataRow = 2
$lastDataRow = 31
$dateColumn = "BF"
$correctedDate = "2007-12-24"

$targetRange = $ws.Range($dateColumn + $firstDataRow + ":" + $dateColumn + $lastDataRow)

# Ensure the cells keep storing a literal text value instead of being
# auto-converted into a date serial by Excel's input parser.
$targetRange.NumberFormat = "@"

foreach ($cell in $targetRange.Cells) {
    if ($cell.Value2 -eq "12-24-2007-08") {
        $cell.Value = $correctedDate
    }
}
